# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 261.3203778131603, 1133.036916526867, 1, 1399.299904967153)
    3 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
    4 = @(0.1190320826869504, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 1.672833113781282)
    5 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 0, 5.586269137925634)
    6 = @(0.01293466051926884, 0.04071648406533734, 22.3905356188092, 10.19245300693656, 0, 32.63663977033036)
    7 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 1133.036916526867, 1, 1138.732270421766)
    8 = @(1.455362044514542, 1.655778082260271, 261.3203778131603, 10.19245300693656, 0, 274.6239709468716)
    9 = @(1.455362044514542, 10.34677158129881, 3.537761648806719, 10.19245300693656, 1, 25.53234828155663)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
